# This edit reorders data rows 5-72 of the "Artfynd" sheet.
# The header (row 1) and rows 2-4 stay fixed; rows 5-72 (68 rows, columns A:AY)
# are permuted into a new order. We read the full block of values, build the
# permuted block in memory, and write it back in one shot so every column
# (including ones not mentioned explicitly in any particular row) follows
# its row to the new position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRange = $ws.Range("A5:AY72")
$data = $srcRange.Value2

$rowCount = $data.GetLength(0)   # 68
$colCount = $data.GetLength(1)   # 51 (A..AY)

# For each new row offset (0-based, row 5 = offset 0), the 0-based offset of
# the source row (also relative to row 5) whose entire contents should be
# copied into that position.
$order = @(48,0,1,2,3,4,5,6,7,8,49,50,51,9,10,52,11,12,53,13,54,14,15,16,17,18,19,20,21,22,55,56,57,58,23,24,59,60,25,61,26,27,62,28,63,64,29,30,31,65,32,33,34,35,36,37,38,39,40,41,42,66,43,67,44,45,46,47)

$newData = New-Object 'object[,]' $rowCount,$colCount

for ($i = 0; $i -lt $rowCount; $i++) {
    $srcRow = $order[$i]
    for ($j = 1; $j -le $colCount; $j++) {
        $newData[$i, $j-1] = $data.GetValue($srcRow + 1, $j)
    }
}

$srcRange.Value2 = $newData
